$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Moorings")

# Recover Date (G2) - fixed deployment date based on cruise reports
$ws.Range("G2").Value = Get-Date -Year 2016 -Month 3 -Day 16 -Hour 0 -Minute 0 -Second 0

# Notes (L2) - glider lost
$ws.Range("L2").Value = "Glider lost"

# Column G (Recover Date) now needs to display the new date value - size it
# to fit, like the other date columns on this sheet.
$ws.Columns.Item(7).ColumnWidth = 9.666666666666666

# Update the selected cell to reflect where the editor left off
$ws.Range("I7").Select()
